$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New PackageTrackNum values for rows 2..22 (column C), in order.
$values = @(
    "320017962708",
    "320017962719",
    "320017962741",
    "320017962774",
    "320017962811",
    "320017962833",
    "320017962866",
    "320017962888",
    "320017962936",
    "320017962958",
    "320017962991",
    "320017963016",
    "320017963049",
    "320017963060",
    "320017963093",
    "320017963119",
    "320017963152",
    "320017963174",
    "320017963200",
    "320017963222",
    "320017963255"
)

# Rows (1-based worksheet rows) where column D mirrors column C's PackageTrackNum value
$mirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $val = $values[$i]

    # Force text storage (so the numeric-looking track number is written as a
    # shared string, matching the workbook's existing convention) and then
    # restore the default "Normal" style so no extra number-format style gets
    # attached to the cell.
    $cellC = $ws.Cells.Item($row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $val
    $cellC.Style = "Normal"

    if ($mirrorRows -contains $row) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $val
        $cellD.Style = "Normal"
    }
}
